# Re-generate of the "library-prep-cs" CodeSystem spreadsheet: refreshes the
# publication Date on the Metadata sheet, bumps the concept Count from 1 to 2,
# and appends a new "unknown" concept row on the Concepts sheet.

$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsConcepts = $wb.Worksheets.Item("Concepts")

# A scratch cell used to manufacture text-typed values (="..." formulas
# always evaluate to text) that we can PasteSpecial(Values) elsewhere so the
# destination keeps the "text" cell type instead of Excel's literal-entry
# auto-number-coercion (e.g. "1", "2" typed directly become numbers).
$scratch = $wsMeta.Cells.Item(30, 10)

# --- Metadata sheet: refresh the generation Date ---
$wsMeta.Range("B8").Value = "2024-09-13T20:57:31+00:00"

# --- Metadata sheet: Count goes from 1 to 2 (still stored as text) ---
$scratch.Formula = '="2"'
$scratch.Copy()
$wsMeta.Cells.Item(22, 2).PasteSpecial(-4163)  # xlPasteValues

# --- Concepts sheet: append the "unknown" concept as row 3 ---
# Copy row 2's formatting down to row 3 first so styles (incl. the blank
# Definition cell) line up exactly with the existing data row.
$wsConcepts.Range("A2:D2").Copy()
$wsConcepts.Range("A3:D3").PasteSpecial(-4122)  # xlPasteFormats

$scratch.Formula = '="1"'
$scratch.Copy()
$wsConcepts.Cells.Item(3, 1).PasteSpecial(-4163)  # Level

$scratch.Formula = '="unknown"'
$scratch.Copy()
$wsConcepts.Cells.Item(3, 2).PasteSpecial(-4163)  # Code

$scratch.Formula = '="Unknown"'
$scratch.Copy()
$wsConcepts.Cells.Item(3, 3).PasteSpecial(-4163)  # Display

# Definition (D3) stays blank, matching the formatted-but-empty D2 cell.

$scratch.Clear()
